$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new timesheet row for Sophia Wilhelmi (row 6), matching the
# formatting already used by the rows above it (rows 3-5).
# Copy the formatting (styles/number formats) from row 5 down into row 6
# first, so the new cells reuse the existing styles instead of creating
# new ones.
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = "Sophia Wilhelmi"
$ws.Range("B6").Value = 42670.479166666664
$ws.Range("C6").Value = 42670.5
